# Updates the HotStock_Top20 ranking table on the active sheet to reflect
# the latest data snapshot (financial news / Eastmoney / 10jqka rankings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target contents for A2:C21 (row 1 is the unchanged header row).
$data = @(
    @("闻泰科技", "闻泰科技", "士兰微"),
    @("大有能源", "华天科技", "闻泰科技"),
    @("合肥城建", "海峡股份", "山子高科"),
    @("三花智控", "白银有色", "楚江新材"),
    @("东信和平", "三花智控", "道生天合"),
    @("白银有色", "东方财富", "华天科技"),
    @("海峡股份", "大有能源", "合肥城建"),
    @("华天科技", "合肥城建", "三花智控"),
    @("山子高科", "贵州茅台", "白银有色"),
    @("睿能科技", "东信和平", "大有能源"),
    @("楚江新材", "平潭发展", "海峡股份"),
    @("平潭发展", "睿能科技", "常山北明"),
    @("贵州茅台", "山子高科", "东信和平"),
    @("东方财富", "东贝集团", "安泰科技"),
    @("超讯通信", "国新能源", "三孚股份"),
    @("北方稀土", "山东墨龙", "远大控股"),
    @("国新能源", "和而泰",   "天际股份"),
    @("士兰微",   "海通发展", "睿能科技"),
    @("三孚股份", "安泰集团", "平潭发展"),
    @("安泰集团", "孚日股份", "华建集团")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
